$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -1226
$ws.Range("H5").Value = 570.5833
$ws.Range("I5").Value = 383.1111
$ws.Range("J5").Value = 1133
$ws.Range("K5").Value = 383.1111
$ws.Range("L5").Value = 1133
$ws.Range("M5").Value = -268.1111
$ws.Range("N5").Value = -1363.3333
$ws.Range("H116").Value = 4247.25
$ws.Range("I116").Value = 4212.8335
$ws.Range("K116").Value = 4212.8335
$ws.Range("M116").Value = -770.8334999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2235.625
$ws.Range("I28").Value = 2235.625
$ws.Range("K28").Value = 2235.625
$ws.Range("M28").Value = -2043.625
$ws.Range("H61").Value = 3681
$ws.Range("I61").Value = 3681
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3681
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3469
$ws.Range("N61").Value = $null
$ws.Range("H99").Value = 2235.625
$ws.Range("I99").Value = 2235.625
$ws.Range("K99").Value = 2235.625
$ws.Range("M99").Value = 759.375
$ws.Range("H136").Value = 3681
$ws.Range("I136").Value = 3681
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11043
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8493
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1014.1
$ws.Range("J20").Value = 1467.3334
$ws.Range("L20").Value = 1467.3334
$ws.Range("N20").Value = -1961.3334
$ws.Range("H54").Value = 4000
$ws.Range("I54").Value = 4000
$ws.Range("K54").Value = 4000
$ws.Range("M54").Value = -3516
$ws.Range("H99").Value = 3337.923
$ws.Range("J99").Value = 2832.6667
$ws.Range("L99").Value = 2832.6667
$ws.Range("N99").Value = -5828.6667
$ws.Range("H110").Value = 102333.336
$ws.Range("I110").Value = 27000
$ws.Range("K110").Value = 27000
$ws.Range("M110").Value = -22910
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 123.30769
$ws.Range("I7").Value = 278
$ws.Range("J7").Value = 26.625
$ws.Range("K7").Value = 278
$ws.Range("L7").Value = 26.625
$ws.Range("M7").Value = -165
$ws.Range("N7").Value = -252.625
$ws.Range("H22").Value = 830
$ws.Range("I22").Value = 830
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 830
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -480
$ws.Range("N22").Value = $null
$ws.Range("H38").Value = 22666.5
$ws.Range("J38").Value = 21749.75
$ws.Range("L38").Value = 21749.75
$ws.Range("N38").Value = -22503.75
$ws.Range("H46").Value = 22666.5
$ws.Range("J46").Value = 21749.75
$ws.Range("L46").Value = 21749.75
$ws.Range("N46").Value = -22171.75
$ws.Range("H88").Value = 24269.084
$ws.Range("J88").Value = 24269.084
$ws.Range("L88").Value = 24269.084
$ws.Range("N88").Value = -25081.084
$ws.Range("H91").Value = 24269.084
$ws.Range("J91").Value = 24269.084
$ws.Range("L91").Value = 24269.084
$ws.Range("N91").Value = -27077.084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2001.8823
$ws.Range("J22").Value = 2001.8823
$ws.Range("L22").Value = 6005.6469
$ws.Range("N22").Value = -6343.6469
$ws.Range("H27").Value = 2001.8823
$ws.Range("J27").Value = 2001.8823
$ws.Range("L27").Value = 6005.6469
$ws.Range("N27").Value = -6209.6469
$ws.Range("H46").Value = 2312.3635
$ws.Range("I46").Value = 472.6
$ws.Range("J46").Value = 3845.5
$ws.Range("K46").Value = 1417.8
$ws.Range("L46").Value = 11536.5
$ws.Range("M46").Value = -1326.8
$ws.Range("N46").Value = -11718.5
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("H80").Value = 4265.3335
$ws.Range("I80").Value = 5898
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 17694
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -16758
$ws.Range("N80").Value = -4872
$ws.Range("H83").Value = 4265.3335
$ws.Range("I83").Value = 5898
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 53082
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -48402
$ws.Range("N83").Value = -18360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("H40").Value = 12000
$ws.Range("I40").Value = 12000
$ws.Range("K40").Value = 12000
$ws.Range("M40").Value = -11849
$ws.Range("H70").Value = 4239.5
$ws.Range("J70").Value = 4239.5
$ws.Range("L70").Value = 4239.5
$ws.Range("N70").Value = -4779.5
$ws.Range("H73").Value = 4239.5
$ws.Range("J73").Value = 4239.5
$ws.Range("L73").Value = 4239.5
$ws.Range("N73").Value = -6111.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 14794.167
$ws.Range("J69").Value = 14794.167
$ws.Range("L69").Value = 14794.167
$ws.Range("N69").Value = -16292.167
$ws.Range("H72").Value = 14794.167
$ws.Range("J72").Value = 14794.167
$ws.Range("L72").Value = 44382.501
$ws.Range("N72").Value = -51870.501
$ws.Range("H132").Value = 1873.5
$ws.Range("I132").Value = 1873.5
$ws.Range("K132").Value = 5620.5
$ws.Range("M132").Value = -3090.5
